# Insert a new daily price record for Naranja (Agrícola del Norte S.A. de Arica)
# as row 40, pushing all existing rows from 40 onward down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 40 (shifts rows 40..127 down to 41..128)
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new record's data
$ws.Range("A40").Value = 1
$ws.Range("B40").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C40").Value = "Arica y Parinacota"
$ws.Range("D40").Value = 44980
$ws.Range("E40").Value = 15
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100102
$ws.Range("H40").Value = "Cítricos"
$ws.Range("I40").Value = 100102005
$ws.Range("J40").Value = "Naranja"
$ws.Range("K40").Value = "Valencia"
$ws.Range("L40").Value = "Tercera"
$ws.Range("M40").Value = 300
$ws.Range("N40").Value = 1150
$ws.Range("O40").Value = 1200
$ws.Range("P40").Value = 1183
$ws.Range("Q40").Value = "`$/kilo (en caja de 20 kilos)"
$ws.Range("R40").Value = "Región de O'Higgins"
$ws.Range("S40").Value = 1183
$ws.Range("T40").Value = 1
